$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 757.6875
$ws.Range("I28").Value = 794.5
$ws.Range("K28").Value = 794.5
$ws.Range("M28").Value = -309.5

$ws.Range("H53").Value = 200.7
$ws.Range("I53").Value = 211.2
$ws.Range("K53").Value = 211.2
$ws.Range("M53").Value = 425.8

$ws.Range("H98").Value = 30168.334
$ws.Range("I98").Value = 25252.5
$ws.Range("J98").Value = 40000
$ws.Range("K98").Value = 25252.5
$ws.Range("L98").Value = 40000
$ws.Range("M98").Value = -23754.5
$ws.Range("N98").Value = -42996

$ws.Range("H122").Value = 30168.334
$ws.Range("I122").Value = 25252.5
$ws.Range("J122").Value = 40000
$ws.Range("K122").Value = 75757.5
$ws.Range("L122").Value = 120000
$ws.Range("M122").Value = -73307.5
$ws.Range("N122").Value = -124900

$ws.Range("H138").Value = 2557.2
$ws.Range("J138").Value = 3036.8667
$ws.Range("L138").Value = 9110.6001
$ws.Range("N138").Value = -19390.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 296.2143
$ws.Range("I2").Value = 309.3846
$ws.Range("K2").Value = 309.3846
$ws.Range("M2").Value = -196.3846

$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 500
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -732

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H41").Value = 5902.375
$ws.Range("I41").Value = 1888.4286
$ws.Range("J41").Value = 34000
$ws.Range("K41").Value = 1888.4286
$ws.Range("L41").Value = 34000
$ws.Range("M41").Value = -1474.4286
$ws.Range("N41").Value = -34828

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H74").Value = 2521.682
$ws.Range("I74").Value = 2536.15
$ws.Range("J74").Value = 2377
$ws.Range("K74").Value = 2536.15
$ws.Range("L74").Value = 2377
$ws.Range("M74").Value = -1662.15
$ws.Range("N74").Value = -4125

$ws.Range("H77").Value = 2521.682
$ws.Range("I77").Value = 2536.15
$ws.Range("J77").Value = 2377
$ws.Range("K77").Value = 12680.75
$ws.Range("L77").Value = 11885
$ws.Range("M77").Value = -8312.75
$ws.Range("N77").Value = -20621

$ws.Range("H97").Value = 1513.2858
$ws.Range("I97").Value = 1053.2727
$ws.Range("J97").Value = 3200
$ws.Range("K97").Value = 1053.2727
$ws.Range("L97").Value = 3200
$ws.Range("M97").Value = -557.2727
$ws.Range("N97").Value = -4192

$ws.Range("H110").Value = 543.4286
$ws.Range("I110").Value = 463.2143
$ws.Range("J110").Value = 703.8570999999999
$ws.Range("K110").Value = 463.2143
$ws.Range("L110").Value = 703.8570999999999
$ws.Range("M110").Value = 1581.7857
$ws.Range("N110").Value = -4793.8571

$ws.Range("H116").Value = 296.2143
$ws.Range("I116").Value = 309.3846
$ws.Range("K116").Value = 309.3846
$ws.Range("M116").Value = 1984.6154

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 296.2143
$ws.Range("I3").Value = 309.3846
$ws.Range("K3").Value = 309.3846
$ws.Range("M3").Value = -195.3846

$ws.Range("H94").Value = 2233.2778
$ws.Range("I94").Value = 1932.5
$ws.Range("J94").Value = 2834.8333
$ws.Range("K94").Value = 1932.5
$ws.Range("L94").Value = 2834.8333
$ws.Range("M94").Value = -1481.5
$ws.Range("N94").Value = -3736.8333

$ws.Range("H107").Value = 3337.5454
$ws.Range("I107").Value = 3306.0476
$ws.Range("K107").Value = 3306.0476
$ws.Range("M107").Value = -1386.0476

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 498
$ws.Range("I7").Value = 498
$ws.Range("K7").Value = 498
$ws.Range("M7").Value = -385

$ws.Range("H17").Value = 1485.6666
$ws.Range("I17").Value = 2003.5
$ws.Range("J17").Value = 450
$ws.Range("K17").Value = 2003.5
$ws.Range("L17").Value = 450
$ws.Range("M17").Value = -1829.5
$ws.Range("N17").Value = -798

$ws.Range("H25").Value = 1233.3334
$ws.Range("I25").Value = 1300
$ws.Range("J25").Value = 1100
$ws.Range("K25").Value = 1300
$ws.Range("L25").Value = 1100
$ws.Range("M25").Value = -1126
$ws.Range("N25").Value = -1448

$ws.Range("H141").Value = 20000
$ws.Range("I141").Value = 20000
$ws.Range("K141").Value = 20000
$ws.Range("M141").Value = -14820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2397.0857
$ws.Range("J4").Value = 3307.6924
$ws.Range("L4").Value = 9923.0772
$ws.Range("N4").Value = -10147.0772

$ws.Range("H44").Value = 1462.5
$ws.Range("I44").Value = 425
$ws.Range("J44").Value = 2500
$ws.Range("K44").Value = 1275
$ws.Range("L44").Value = 7500
$ws.Range("M44").Value = -877
$ws.Range("N44").Value = -8296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31136

$ws.Range("H80").Value = 2399.8
$ws.Range("J80").Value = 2399.8
$ws.Range("L80").Value = 2399.8
$ws.Range("N80").Value = -4395.8

$ws.Range("H83").Value = 2399.8
$ws.Range("J83").Value = 2399.8
$ws.Range("L83").Value = 11999
$ws.Range("N83").Value = -21983

$ws.Range("H98").Value = 12500
$ws.Range("J98").Value = 12500
$ws.Range("L98").Value = 12500
$ws.Range("N98").Value = -18490

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1701000
$ws.Range("I40").Value = 51500
$ws.Range("K40").Value = 51500
$ws.Range("M40").Value = -51364

$ws.Range("H46").Value = 337995.84
$ws.Range("I46").Value = 668833.3
$ws.Range("J46").Value = 7158.3335
$ws.Range("K46").Value = 668833.3
$ws.Range("L46").Value = 7158.3335
$ws.Range("M46").Value = -668645.3
$ws.Range("N46").Value = -7534.3335

$ws.Range("H68").Value = 1800
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1800
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H104").Value = 35257.5
$ws.Range("J104").Value = 35257.5
$ws.Range("L104").Value = 35257.5
$ws.Range("N104").Value = -42245.5

$ws.Range("H132").Value = 1798
$ws.Range("I132").Value = 1798
$ws.Range("K132").Value = 5394
$ws.Range("M132").Value = -2864
